$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet "Q1_20_21" : SoT/F9 (rows 3-4) -> Mars/Sea of Tranquility/Apollo 13/
#                    Falcon 9/Columbia (rows 3-7)
# ======================================================================
$ws1 = $wb.Worksheets.Item("Q1_20_21")

# Before: row3 = SoT, row4 = F9
# Insert a blank row at 3: row3 = blank(new), row4 = SoT, row5 = F9
$ws1.Rows.Item(3).Insert()

# New row 3: Mars
$ws1.Cells.Item(3, 2).Value = "Mars"
$ws1.Cells.Item(3, 3).Value = 28369
$ws1.Cells.Item(3, 4).Value = 14.58
$ws1.Cells.Item(3, 5).Value = 12.98
$ws1.Cells.Item(3, 6).Value = "Very High"
$ws1.Cells.Item(3, 7).Value = "Very High"
$ws1.Cells.Item(3, 8).Value = "Very High"
$ws1.Cells.Item(3, 9).Value = 2089
$ws1.Cells.Item(3, 10).Value = 30458
$ws1.Cells.Item(3, 11).Value = "All you need is love, love is all you need "

# Row 4 is now the old "SoT" row -> rename to "Sea of Tranquility" (data unchanged)
$ws1.Cells.Item(4, 2).Value = "Sea of Tranquility"

# Before: row5 = F9
# Insert a blank row at 5: row5 = blank(new), row6 = F9
$ws1.Rows.Item(5).Insert()

# New row 5: Apollo 13
$ws1.Cells.Item(5, 2).Value = "Apollo 13"
$ws1.Cells.Item(5, 3).Value = 1985
$ws1.Cells.Item(5, 4).Value = 2.3
$ws1.Cells.Item(5, 5).Value = 2.3
$ws1.Cells.Item(5, 6).Value = "High"
$ws1.Cells.Item(5, 9).Value = 833
$ws1.Cells.Item(5, 10).Value = 3494

# Row 6 is now the old "F9" row -> rename to "Falcon 9" (data unchanged:
# 1356, 1.46, 0.74, Medium, N/A, N/A, 2956, 4312)
$ws1.Cells.Item(6, 2).Value = "Falcon 9"

# Insert a blank row at 7 for the new Columbia entry
$ws1.Rows.Item(7).Insert()

# New row 7: Columbia
$ws1.Cells.Item(7, 2).Value = "Columbia"
$ws1.Cells.Item(7, 4).Value = 0.38
$ws1.Cells.Item(7, 5).Value = 0.63
$ws1.Cells.Item(7, 6).Value = "Poor"
$ws1.Cells.Item(7, 9).Value = 1172
$ws1.Cells.Item(7, 10).Value = 738.36

# ======================================================================
# Sheet "Q4_19_20" : SoT/A13/Columbia (rows 3-5) -> Mars/Sea of Tranquility/
#                    Apollo 11/Apollo 13/Falcon 9/Columbia (rows 3-8)
# ======================================================================
$ws2 = $wb.Worksheets.Item("Q4_19_20")

# Before: row3 = SoT, row4 = A13, row5 = Columbia
# Insert a blank row at 3: row3 = blank(new), row4 = SoT, row5 = A13, row6 = Columbia
$ws2.Rows.Item(3).Insert()

# New row 3: Mars
$ws2.Cells.Item(3, 2).Value = "Mars"
$ws2.Cells.Item(3, 3).Value = 30292.2
$ws2.Cells.Item(3, 4).Value = 21.45
$ws2.Cells.Item(3, 5).Value = 19.72
$ws2.Cells.Item(3, 6).Value = "Very High"
$ws2.Cells.Item(3, 9).Value = 1481.6
$ws2.Cells.Item(3, 10).Value = 31773.8
$ws2.Cells.Item(3, 11).Value = "Hello is it me you’re looking for"

# Row 4 is now the old "SoT" row -> rename to "Sea of Tranquility" (data + K4
# narrative unchanged)
$ws2.Cells.Item(4, 2).Value = "Sea of Tranquility"

# Before: row5 = A13, row6 = Columbia
# Insert a blank row at 5: row5 = blank(new), row6 = A13, row7 = Columbia
$ws2.Rows.Item(5).Insert()

# New row 5: Apollo 11 (only the name is populated)
$ws2.Cells.Item(5, 2).Value = "Apollo 11"

# Row 6 is now the old "A13" row -> rename to "Apollo 13" (data unchanged)
$ws2.Cells.Item(6, 2).Value = "Apollo 13"

# Before: row7 = Columbia
# Insert a blank row at 7: row7 = blank(new), row8 = Columbia
$ws2.Rows.Item(7).Insert()

# New row 7: Falcon 9
$ws2.Cells.Item(7, 2).Value = "Falcon 9"
$ws2.Cells.Item(7, 3).Value = 2952
$ws2.Cells.Item(7, 4).Value = 1.54
$ws2.Cells.Item(7, 5).Value = 0.78
$ws2.Cells.Item(7, 6).Value = "Medium"
$ws2.Cells.Item(7, 9).Value = 2831
$ws2.Cells.Item(7, 10).Value = 4364

# Row 8 is the old "Columbia" row; its data is unchanged
# (0.38, 0.63, Poor, 1172, 738.36)
